$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Add / Edit / Delete" section paragraph: insert "'s methods" so that
#    "...with help from the model class, so that..." becomes
#    "...with help from the model class's methods, so that..."
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(55)
$rng1 = $p1.Range
$null = $rng1.Find.Execute(
    "with help from the model class, so that",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "with help from the model class" + [char]0x2019 + "s methods, so that",
    2)

# ---------------------------------------------------------------------------
# 2) "Export" section paragraph: merge "This feat" + "ure, implemented in "
#    into a single run reading "This feature, implemented in " (text itself
#    is unchanged, this just normalizes the run split).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(59)
$rng2 = $p2.Range
$null = $rng2.Find.Execute(
    "This feature, implemented in ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This feature, implemented in ",
    2)

# ---------------------------------------------------------------------------
# 3) Same "Export" paragraph: merge " format." + " " into a single run
#    reading " format. " (adds a trailing space after the period).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(59)
$rng3 = $p3.Range
$null = $rng3.Find.Execute(
    " format. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " format. ",
    2)

# ---------------------------------------------------------------------------
# 4) "View IEEE formatted citation" section paragraph: add " in code" right
#    before the trailing period, and relocate the "_GoBack" bookmark so it
#    sits between "...stylesheet in code" and the final ".".
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(64)
$rng4 = $p4.Range
$null = $rng4.Find.Execute(
    "referencing the supported XML stylesheet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "referencing the supported XML stylesheet in code.",
    2)

if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
}

$p4b = $d.Paragraphs(64)
$rng4b = $p4b.Range
$bookmarkPos = $rng4b.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
